$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.996.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.54%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.262.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.70%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'253.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.12%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.84%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'71.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.05%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.675"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +16.41%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.10%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'39.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.64%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0980"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.22%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'59.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.58%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'7.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.44%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.14%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.599.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.78%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.887"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.28%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'14.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.19%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.260.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.23%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'42.922.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.56%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.0₃0983"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.62%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.03%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'73.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.01%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'238.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.21%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.67%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +0.87%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'11.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.49%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.22%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -2.51%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -1.33%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +8.24%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'168.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.12%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D33").Value = "'6.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +14.88%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +8.35%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0774"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.78%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +1.67%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'29.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +11.34%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +1.97%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -1.09%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.0323"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +7.27%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +3.44%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +2.28%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'12.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.53%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'64.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.29%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'5.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.40%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.203"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.26%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +2.12%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.06%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -4.91%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -0.17%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +1.08%  "
$ws.Range("E51").Style = "Normal"
